# Weekly update: a new week's price row is inserted at row 168 (above the
# current top data row), which pushes all the existing data rows (168-200)
# down by one (to 169-201). The newly inserted row 168 is then populated
# with this week's values; the "constant" descriptive columns are copied
# forward from the row that used to occupy that slot, while the Date (D)
# and Volumen (J) columns get the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 168; rows 168:200 shift down to 169:201
# (formatting/styles of the row, e.g. the date style on column D, move
# along with it, matching native Excel "Insert" behaviour).
$ws.Rows(168).Insert()

# Populate the newly inserted row 168 with the new week's data.
$ws.Cells.Item(168, 1).Value = 8
$ws.Cells.Item(168, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(168, 3).Value = "Coquimbo"
$ws.Cells.Item(168, 4).Value = 44476
$ws.Cells.Item(168, 5).Value = 4
$ws.Cells.Item(168, 6).Value = 100114013
$ws.Cells.Item(168, 7).Value = "Zanahoria"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 560
$ws.Cells.Item(168, 11).Value = 6000
$ws.Cells.Item(168, 12).Value = 7000
$ws.Cells.Item(168, 13).Value = 6500
$ws.Cells.Item(168, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(168, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(168, 16).Value = 325
$ws.Cells.Item(168, 17).Value = 20
$ws.Cells.Item(168, 18).Value = "Hortaliza"
